$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Shift month labels forward by one month
$ws.Range("C1").Value = "abril"
$ws.Range("D1").Value = "mayo"
$ws.Range("E1").Value = "junio"
$ws.Range("F1").Value = "julio"

# Swap column widths for D and E to match new header content.
# (ColumnWidth set here accounts for the stored/display width offset
# so the saved OOXML "width" attribute ends up as 10 and 11.)
$ws.Columns.Item(4).ColumnWidth = 9.166666666666666
$ws.Columns.Item(5).ColumnWidth = 10.166666666666666
